$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the Heading1
#    title paragraph ("Play American Dad Free Slot Game Review").
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$d = $word.ActiveDocument
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$fullMetaText = "Meta description: Play American Dad for free and read our review of the slot game inspired by the television show. Enjoy a wide range of bonus features and chance to win big."
$metaRange.Text = $fullMetaText

$boldStart = $metaRange.Start
$boldEnd = $boldStart + ("Meta description").Length
$boldRange = $d.Range($boldStart, $boldEnd)
$boldRange.Bold = 1

# Leave a leading empty run, matching the document's authoring convention.
$metaPara = $d.Paragraphs.Item(2)
$leadStart = $metaPara.Range.Start
$leadRange = $d.Range($leadStart, $leadStart)
$leadRange.InsertBefore("")

# ---------------------------------------------------------------------------
# 2) Remove the duplicated bold "Play American Dad Free Slot Game Review"
#    paragraph near the end of the document, and rewrite the italic
#    paragraph's text with the new image-generation prompt.
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$dupTitlePara.Range.Delete()

$d = $word.ActiveDocument
$count = $d.Paragraphs.Count
$italicPara = $d.Paragraphs.Item($count)
$italicRange = $italicPara.Range
$textOnly = $d.Range($italicRange.Start, $italicRange.End - 1)
$newImageText = "Create a feature image for American Dad that captures the cartoon style of the game and features a happy Maya warrior with glasses. The image should include the iconic characters from the American Dad TV show, including Stan, Roger, Francine, Hayley, Wheels and Legman. The Maya warrior should be positioned as a playful addition to the group, perhaps holding a slot machine lever or standing next to a pile of coins. The background should showcase Langley Falls and the American flag. The overall tone should be fun and lighthearted, capturing the essence of the game's zany and amusing adventures."
$textOnly.Text = $newImageText
